# Fill rows 2-20 with the same "fn_val / ln_val / address" pattern that was
# already present in row 1, extending the used range of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = "fn_val"
    $ws.Cells.Item($r, 2).Value = "ln_val"
    $ws.Cells.Item($r, 3).Value = "address"
}

# Move/restore the active selection to match the author's final cursor
# position after filling in the new rows.
$ws.Range("F19").Select()
